$wb = $excel.ActiveWorkbook

# --- Sheet "Debts": insert a new "active" column at the front ---
$wsDebts = $wb.Worksheets.Item("Debts")
$wsDebts.Columns.Item(1).Insert()
$wsDebts.Range("A1").Value = "active"

# --- Sheet "Fixed Assets": insert a new "active" column at the front ---
$wsFA = $wb.Worksheets.Item("Fixed Assets")
$wsFA.Columns.Item(1).Insert()
$wsFA.Range("A1").Value = "active"

# --- Selections / active sheet state ---
$wsDebts.Range("C11").Select()
$wsFA.Range("C12").Select()
$wsFA.Activate()
